$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15: Potentiometre 10K line item ---
# Unit price drops from 0.51 to 0.37
$ws.Range("E15").Value = 0.37

# Total-price formula's flat add-on changes from +2 to +0.99 (new total 2.1)
$ws.Range("F15").Formula = "=E15*D15+0.99"

# --- Update the AliExpress product link text shown in H15 and its target ---
$newUrl = "https://fr.aliexpress.com/item/1005002766893077.html?spm=a2g0o.detail.1000014.28.6d93d6c1fPOoi0&gps-id=pcDetailBottomMoreOtherSeller&scm=1007.40050.281175.0&scm_id=1007.40050.281175.0&scm-url=1007.40050.281175.0&pvid=c4eb49e9-017c-42e1-bc5a-36d0b91caa96&_t=gps-id:pcDetailBottomMoreOtherSeller,scm-url:1007.40050.281175.0,pvid:c4eb49e9-017c-42e1-bc5a-36d0b91caa96,tpp_buckets:668%232846%238116%232002&pdp_ext_f=%7B%22sku_id%22%3A%2212000022084500624%22%2C%22sceneId%22%3A%2230050%22%7D&pdp_npi=2%40dis%21EUR%210.42%210.37%21%21%21%21%21%402101f6b416658300092574280e1af7%2112000022084500624%21rec"

# Update the visible cell text (keeps the existing cell style/format intact)
$ws.Range("H15").Value = $newUrl

# Best-effort: repoint the existing hyperlink on H15 to the new address as well,
# without disturbing any of the other hyperlinks on the sheet.
for ($i = 1; $i -le $ws.Hyperlinks.Count; $i++) {
    $hl = $ws.Hyperlinks.Item($i)
    if ($hl.Range.Address() -eq "`$H`$15") {
        $hl.Address = $newUrl
    }
}

# --- Selection cursor moved from I15 to I25 (cosmetic, matches sheetView) ---
$ws.Range("I25").Select()
